$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Insert a new row for "ffmpeg" right before the Firefox row (row 15).
# ---------------------------------------------------------------------------
$ws.Rows.Item(15).Insert() | Out-Null
$ws.Range("A15").Value = "ffmpeg"
$ws.Range("B15").Value = "FFmpeg"
$ws.Range("C15").Value = "2020-09-27 12:31"
$ws.Range("D15").Value = "https://github.com/BtbN/FFmpeg-Builds/releases"
$ws.Range("E15").Value = "https://github.com/BtbN/FFmpeg-Builds/releases"
$ws.Range("F15").Value = '//div[@class="f1 flex-auto min-width-0 text-normal"]/a/text()'

# ---------------------------------------------------------------------------
# 2. Insert a new row for "spotify" right before the Sublime Text row.
#    After step 1, Sublime Text moved from row 27 to row 28.
# ---------------------------------------------------------------------------
$ws.Rows.Item(28).Insert() | Out-Null
$ws.Range("A28").Value = "spotify"
$ws.Range("B28").Value = "Spotify"
$ws.Range("C28").Value = "1.142.622.0"
$ws.Range("D28").Value = "https://en.wikipedia.org/wiki/Spotify"
$ws.Range("E28").Value = "http://download.spotify.com/SpotifyFullSetup.exe"
$ws.Range("F28").Value = '//table[@class="infobox vevent"]/tbody/tr/td//text()'

# ---------------------------------------------------------------------------
# 3. Delete the "yt_music" (YouTube Music) row.
#    After steps 1-2, it sits at row 32.
# ---------------------------------------------------------------------------
$ws.Rows.Item(32).Delete() | Out-Null

# ---------------------------------------------------------------------------
# 4. Update version numbers / values that changed in rows that kept their
#    relative position (rows 2-14, unaffected by the inserts above).
# ---------------------------------------------------------------------------
$ws.Range("C2").Value = "0.37.0"              # Visual C Redistributable
$ws.Range("C5").Value = "5.0.1"               # Calibre
$ws.Range("C6").Value = "1.3.16 2020/07/29"   # cmder
$ws.Range("C7").Value = "8.8.8"               # CrystalDiskInfo
$ws.Range("C9").Value = "106.4.368"           # Dropbox
$ws.Range("C11").Value = "2020-09-09"         # Peace Equalizer
$ws.Range("C12").Value = "1.4.1.992"          # Search Everything
$ws.Range("C13").Value = "12.06"              # Exiftool

# ---------------------------------------------------------------------------
# 5. Update version numbers / values for rows 16-27 (shifted down by 1 row
#    because of the ffmpeg insert in step 1).
# ---------------------------------------------------------------------------
$ws.Range("C16").Value = "81.0"               # Firefox
$ws.Range("C17").Value = "2.28.0"             # Git
$ws.Range("C20").Value = "6.30"               # HWInfo
$ws.Range("C21").Value = "12.10.9.3"          # iTunes
$ws.Range("C22").Value = "15.7.0"             # K-Lite Codec
$ws.Range("C23").Value = "50.0.0"             # MKVToolnix
$ws.Range("C25").Value = "4.4.152"            # Open Shell

# ---------------------------------------------------------------------------
# 6. Update version numbers / values for rows 30-32 (shifted down by 2 rows
#    because of the ffmpeg + spotify inserts, then back up by 1 because of
#    the yt_music delete).
# ---------------------------------------------------------------------------
$ws.Range("C30").Value = "1.49"               # Visual Studio Code
$ws.Range("C31").Value = "5.17.7"             # WinSCP
$ws.Range("C32").Value = "15"                 # Java SE

# ---------------------------------------------------------------------------
# 7. Fix up sheet selection / active cell (cosmetic, matches target diff).
# ---------------------------------------------------------------------------
$ws.Range("A1").Select() | Out-Null

Write-Host "Edit complete"
